$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8930620761606677
$ws.Range("C2").Value = 0.9823589743589743

$ws.Range("B3").Value = 0.7219613980177361
$ws.Range("C3").Value = 0.7292307692307692

$ws.Range("B4").Value = 0.7219613980177361
$ws.Range("C4").Value = 0.7292307692307692

$ws.Range("B5").Value = 0.782472613458529
$ws.Range("C5").Value = 0.762051282051282

$ws.Range("B6").Value = 0.7219613980177361
$ws.Range("C6").Value = 0.7292307692307692

$ws.Range("B7").Value = 0.8085550339071466
$ws.Range("C7").Value = 0.7965128205128206

$ws.Range("B8").Value = 0.7219613980177361
$ws.Range("C8").Value = 0.7292307692307692

$ws.Range("B9").Value = 0.7219613980177361
$ws.Range("C9").Value = 0.7292307692307692

$ws.Range("B10").Value = 0.7219613980177361
$ws.Range("C10").Value = 0.7292307692307692

$ws.Range("B11").Value = 0.7219613980177361
$ws.Range("C11").Value = 0.7292307692307692

$ws.Range("B12").Value = 0.8038601982263954
$ws.Range("C12").Value = 0.7669743589743589
